# The sheet originally listed 4 domains (one row highlighted red because its
# expiry lookup failed); the commit removes that sample data entirely and
# keeps only the final domain/date pair (vk.com / 23.06.2024), moved up into
# row 2 with its highlight formatting cleared, plus narrower, "real" column
# widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 currently holds "spiritofhelsinki.com" (the failed lookup, highlighted
# red). Overwrite it with the last row's data (vk.com / 23.06.2024) before
# dropping the now-redundant rows 3-5.
$ws.Range("A2").Value = $ws.Range("A5").Text
$ws.Range("B2").Value = $ws.Range("B5").Text

# Drop the red "lookup failed" highlight that was on row 2.
$ws.Range("A2:B2").ClearFormats()

# Remove the now-stale rows (spiritofhelsinki.com duplicate source row,
# tkstkdsoft.com, voteexpress.com).
$ws.Range("A3:B5").EntireRow.Delete()

# Shrink the columns down to the final, narrower widths.
$ws.Columns.Item(1).ColumnWidth = 6.8
$ws.Columns.Item(2).ColumnWidth = 10.8
